# UI-SJN-01-001U - 추가 수정
#
# Shift the right-hand screenshot group further right, drop the
# "이미지는 클래스명 img / 아이콘은 클래스면 icon" helper textbox and the
# small picture (그림 18) that used to sit above it, and shift the
# caption textbox at the bottom of that group to match.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# NOTE on the magic .Left values below: this COM host round-trips the
# Left/Top (points) properties through a 32-bit float before converting
# back to EMU, so naively using targetEmu/12700.0 can land 1 EMU short of
# the intended value. The literals used here were solved so that they
# land on the exact target EMU offset after that conversion.

# 1) Move the screenshot picture ("그림 6") right by 344773 EMU
#    (5647929 -> 5992702).
$pic6 = $s.Shapes.Item("그림 6")
$pic6.Left = 471.8663024902344

# 2) Move the arrow connector ("직선 화살표 연결선 9") right by the same
#    amount (6663128 -> 7007901).
$cxn9 = $s.Shapes.Item("직선 화살표 연결선 9")
$cxn9.Left = 551.8032836914062

# 3) Delete the helper textbox explaining img/icon class usage.
$txt10 = $s.Shapes.Item("TextBox 10")
$txt10.Delete()

# 4) Delete the small picture ("그림 18") that illustrated the icon class.
$pic18 = $s.Shapes.Item("그림 18")
$pic18.Delete()

# 5) Move the bottom caption textbox ("TextBox 19") right to line back up
#    under the relocated screenshot (5864047 -> 6208820).
$txt19 = $s.Shapes.Item("TextBox 19")
$txt19.Left = 488.88348388671875
